$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1. Update the time_taken (column F) values on the "data" sheet.
# ---------------------------------------------------------------------------
$newTimes = @(
    "2021-10-05 14:34:57.892950",
    "2021-10-05 14:34:57.892958",
    "2021-10-05 14:34:57.892961",
    "2021-10-05 14:34:57.892963",
    "2021-10-05 14:34:57.892966",
    "2021-10-05 14:34:57.892969",
    "2021-10-05 14:34:57.892972",
    "2021-10-05 14:34:57.892974",
    "2021-10-05 14:34:57.892977",
    "2021-10-05 14:34:57.892980",
    "2021-10-05 14:34:57.892982",
    "2021-10-05 14:34:57.892985",
    "2021-10-05 14:34:57.892987",
    "2021-10-05 14:34:57.892989",
    "2021-10-05 14:34:57.892992",
    "2021-10-05 14:34:57.892994",
    "2021-10-05 14:34:57.892997",
    "2021-10-05 14:34:57.893000",
    "2021-10-05 14:34:57.893002",
    "2021-10-05 14:34:57.893005",
    "2021-10-05 14:34:57.893007",
    "2021-10-05 14:34:57.893009",
    "2021-10-05 14:34:57.893012",
    "2021-10-05 14:34:57.893014",
    "2021-10-05 14:34:57.893017",
    "2021-10-05 14:34:57.893020"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Range("F$row").Value = $newTimes[$i]
}

# ---------------------------------------------------------------------------
# 2. Add a new "metadata" sheet after the "data" sheet.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row (reuse the bold/centered/bordered header style from "data").
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Data row.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Multiple pterygium syndrome_Fetal akinesia sequence"
$ws.Range("C2").Value = 139
# "1.0" must stay a text value (not be coerced into the number 1), so force
# text entry with a leading apostrophe and then drop the resulting
# quote-prefix style so the cell keeps the default (unstyled) formatting.
$ws.Range("D2").Value = "'1.0"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "2021-06-14T08:25:44.029926Z"
$ws.Range("F2").Value = "2021-10-05 14:34:57.889323"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/139/?format=json"

$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$dataSheet.Select()
